$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.076.90"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.832.53"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").Value = "0.9965"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "241.94"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "0.6171"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.07463"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.2930"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "23.07"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "0.07659"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "1.827.02"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "5.001"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "0.6747"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "'82.90"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "'0.000009199"
$ws.Range("E16").Value = "  -4.15%  "
$ws.Range("D17").Value = "5.905"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "29.046.92"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "2.075.52"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "239.61"
$ws.Range("E20").Value = "  +6.14%  "
$ws.Range("D21").Value = "'12.70"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "7.203"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("D24").Value = "0.9982"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "158.94"
$ws.Range("D26").Value = "0.1406"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "8.502"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "17.89"
$ws.Range("D29").Value = "1.497"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "0.05594"
$ws.Range("E30").Value = "  +3.06%  "
$ws.Range("D31").Value = "4.141"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").Value = "4.119"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").Value = "1.201"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").Value = "1.845"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").Value = "0.7416"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").Value = "1.142"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("D37").Value = "2.655"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "0.01786"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").Value = "1.214.67"
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("D41").Value = "6.422"
$ws.Range("E41").Value = "  -3.30%  "
$ws.Range("D42").Value = "0.8958"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").Value = "0.9985"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "101.39"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "1.974.04"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "65.47"
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").Value = "0.5081"
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("D49").Value = "0.4066"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "9.179"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("D51").Value = "0.05806"
$ws.Range("E51").Value = "  +0.45%  "
